# Update the personal-data sheet: replace the generic "id1/PersonName" example
# row with a real person's record (Juan Torres Pardo) and extend the table
# with the additional fields used by the new form (birth date, postal
# address, nationality, DNI/NIF and polling station).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had two extra, empty formatted rows (3 and 4) that are no
# longer needed - get rid of them first.
$ws.Rows("3:4").Delete() | Out-Null

# --- Row 1: column headers -------------------------------------------------
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Apellidos"
$ws.Range("C1").Value = "Correo electrónico"
$ws.Range("D1").Value = "Fecha nacimiento"
$ws.Range("E1").Value = "Dirección postal"
$ws.Range("F1").Value = "Nacionalidad"
$ws.Range("G1").Value = "DNI"
$ws.Range("H1").Value = "NIF"
$ws.Range("I1").Value = "pollingStation"

# --- Row 2: the actual record ----------------------------------------------
$ws.Range("A2").Value = "Juan"
$ws.Range("B2").Value = "Torres Pardo"

# C2 keeps the existing hyperlink formatting/relationship; only its
# displayed text changes.
$ws.Range("C2").Value = "juan@example.com"

# Birth date, stored as a real date serial with a short-date number format.
$ws.Range("D2").Value = 31330
$ws.Range("D2").NumberFormat = "mm-dd-yy"

$ws.Range("E2").Value = "C/ Federico García Lorca 2"
$ws.Range("F2").Value = "Español"
# G2 (DNI) is intentionally left blank.
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1

# Match the active selection left behind in the edited workbook.
$ws.Range("G2").Select() | Out-Null
